# Insert a new weekly price row for "Vega Modelo de Temuco" / Locoto.
# The new observation (2023-01-05, volumen 140, precio 2500) is inserted
# at row 54, pushing the existing rows 54-58 down to 55-59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 54-58 down by one to make room for the new row.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly record.
$ws.Range("A54").Value = 10
$ws.Range("B54").Value = "Vega Modelo de Temuco"
$ws.Range("C54").Value = "La Araucanía"
$ws.Range("D54").Value = 44931
$ws.Range("E54").Value = 9
$ws.Range("F54").Value = 100112042
$ws.Range("G54").Value = "Locoto"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 140
$ws.Range("K54").Value = 2500
$ws.Range("L54").Value = 2500
$ws.Range("M54").Value = 2500
$ws.Range("N54").Value = "$/kilo"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 2500
$ws.Range("Q54").Value = 1
$ws.Range("R54").Value = "Hortaliza"
